# Update SnippetID (column H) values on the "Voice Lines - main" sheet to
# reflect the refactored snippet grouping. Rows that share a SnippetID
# (e.g. multiple lines belonging to the same snippet) get the same new
# value, matching the original shared-string reuse.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "uXOa"
    3  = "uXOa"
    4  = "N3Ke"
    5  = "9Xte"
    6  = "Bpuz"
    7  = "8KkF"
    8  = "igag"
    9  = "uQLC"
    10 = "5N0e"
    11 = "I2Ri"
    12 = "7cC0"
    13 = "d5aj"
    14 = "AJbt"
    15 = "QxNF"
    16 = "uR5r"
    17 = "zyiY"
    18 = "zyiY"
    19 = "zyiY"
    20 = "zyiY"
    21 = "zyiY"
    22 = "1gvn"
    23 = "pnRQ"
    24 = "arGM"
    25 = "m7xl"
    26 = "ky8X"
    27 = "ky8X"
    28 = "roOK"
    29 = "oUh7"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 8).Value = $updates[$row]
}
